$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add new "venue" column header (C1), matching the existing
#     header style used by A1/B1 ---
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C1").Value = "venue"

# --- Row 2: update existing entry, add venue ---
$ws.Range("A2").Value = "2023-11-29 14:45:20+08:00"
$ws.Range("B2").Value = "Student - Lim Qin Xin - Startup Springboard"
$ws.Range("C2").Value = "Microsoft Teams Meeting"

# --- Row 3: update existing entry, add venue ---
$ws.Range("A3").Value = "2023-11-29 14:39:47+08:00"
$ws.Range("B3").Value = "Week 12 Cohort 2 In-class Annotation"
$ws.Range("C3").Value = "NA"

# --- Row 4: brand new row ---
$ws.Range("A4").Value = "2023-11-29 14:23:00+08:00"
$ws.Range("B4").Value = "arrangement for MA in Week 13 and 14"
$ws.Range("C4").Value = "NA"

# --- Row 5: brand new row, venue column left blank (present but empty) ---
$ws.Range("A5").Value = "2023-11-29 11:48:36+08:00"
$ws.Range("B5").Value = "Something's Brewing! Accenture University Innovation Challenge 2024"
$ws.Range("C5").Formula = "="""""
